$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5 (record ID=4): update Total_Due and Last_Payment_Date
$ws.Range("D5").Value = 22050.2205
$ws.Range("E5").Value = "2025-03-27 17:50:33"

# Row 6 (record ID=5): update Total_Due and Last_Payment_Date
$ws.Range("D6").Value = 7840.0784
$ws.Range("E6").Value = "2025-03-27 17:47:33"
